$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Tnfsf18"
$ws.Cells.Item(2, 3).Value = "Tnfrsf18"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.1244153333333333
$ws.Cells.Item(2, 8).Value = 0.373246
$ws.Cells.Item(2, 9).Value = 0.1599094129708596
$ws.Cells.Item(2, 10).Value = 0.1599094129708596
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.387697
$ws.Cells.Item(2, 14).Value = 1.163091
$ws.Cells.Item(2, 15).Value = 0.09541986830791312
$ws.Cells.Item(2, 16).Value = 0.09541986830791316
$ws.Cells.Item(2, 17).Value = 0.04823545148733334
$ws.Cells.Item(2, 18).Value = 0.434119063386
$ws.Cells.Item(2, 19).Value = 0.01525853512687511
$ws.Cells.Item(2, 20).Value = 0.01525853512687512

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Tnfsf18"
$ws.Cells.Item(3, 3).Value = "Tnfrsf18"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.1244153333333333
$ws.Cells.Item(3, 8).Value = 0.373246
$ws.Cells.Item(3, 9).Value = 0.1599094129708596
$ws.Cells.Item(3, 10).Value = 0.1599094129708596
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 2.236532666666667
$ws.Cells.Item(3, 14).Value = 6.709598
$ws.Cells.Item(3, 15).Value = 0.5504547430588297
$ws.Cells.Item(3, 16).Value = 0.5504547430588298
$ws.Cells.Item(3, 17).Value = 0.2782589572342222
$ws.Cells.Item(3, 18).Value = 2.504330615108
$ws.Cells.Item(3, 19).Value = 0.08802289482956278
$ws.Cells.Item(3, 20).Value = 0.0880228948295628

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Tnfsf18"
$ws.Cells.Item(4, 3).Value = "Tnfrsf18"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.1244153333333333
$ws.Cells.Item(4, 8).Value = 0.373246
$ws.Cells.Item(4, 9).Value = 0.1599094129708596
$ws.Cells.Item(4, 10).Value = 0.1599094129708596
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.998142
$ws.Cells.Item(4, 14).Value = 2.994426
$ws.Cells.Item(4, 15).Value = 0.2456624069636779
$ws.Cells.Item(4, 16).Value = 0.2456624069636779
$ws.Cells.Item(4, 17).Value = 0.124184169644
$ws.Cells.Item(4, 18).Value = 1.117657526796
$ws.Cells.Item(4, 19).Value = 0.03928373128657013
$ws.Cells.Item(4, 20).Value = 0.03928373128657014

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Tnfsf18"
$ws.Cells.Item(5, 3).Value = "Tnfrsf18"
$ws.Cells.Item(5, 4).Value = "Resolving-Mac"
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0.3333333333333333
$ws.Cells.Item(5, 7).Value = 0.1244153333333333
$ws.Cells.Item(5, 8).Value = 0.373246
$ws.Cells.Item(5, 9).Value = 0.1599094129708596
$ws.Cells.Item(5, 10).Value = 0.1599094129708596
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.4406919999999999
$ws.Cells.Item(5, 14).Value = 1.322076
$ws.Cells.Item(5, 15).Value = 0.1084629816695792
$ws.Cells.Item(5, 16).Value = 0.1084629816695792
$ws.Cells.Item(5, 17).Value = 0.05482884207733332
$ws.Cells.Item(5, 18).Value = 0.4934595786959999
$ws.Cells.Item(5, 19).Value = 0.01734425172785151
$ws.Cells.Item(5, 20).Value = 0.01734425172785151

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Tnfsf18"
$ws.Cells.Item(6, 3).Value = "Tnfrsf18"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 2
$ws.Cells.Item(6, 6).Value = 0.6666666666666666
$ws.Cells.Item(6, 7).Value = 0.2127573333333333
$ws.Cells.Item(6, 8).Value = 0.638272
$ws.Cells.Item(6, 9).Value = 0.2734542388551691
$ws.Cells.Item(6, 10).Value = 0.2734542388551692
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 0.387697
$ws.Cells.Item(6, 14).Value = 1.163091
$ws.Cells.Item(6, 15).Value = 0.09541986830791312
$ws.Cells.Item(6, 16).Value = 0.09541986830791316
$ws.Cells.Item(6, 17).Value = 0.08248537986133334
$ws.Cells.Item(6, 18).Value = 0.742368418752
$ws.Cells.Item(6, 19).Value = 0.02609296745980086
$ws.Cells.Item(6, 20).Value = 0.02609296745980087

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Tnfsf18"
$ws.Cells.Item(7, 3).Value = "Tnfrsf18"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 2
$ws.Cells.Item(7, 6).Value = 0.6666666666666666
$ws.Cells.Item(7, 7).Value = 0.2127573333333333
$ws.Cells.Item(7, 8).Value = 0.638272
$ws.Cells.Item(7, 9).Value = 0.2734542388551691
$ws.Cells.Item(7, 10).Value = 0.2734542388551692
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 2.236532666666667
$ws.Cells.Item(7, 14).Value = 6.709598
$ws.Cells.Item(7, 15).Value = 0.5504547430588297
$ws.Cells.Item(7, 16).Value = 0.5504547430588298
$ws.Cells.Item(7, 17).Value = 0.4758387260728889
$ws.Cells.Item(7, 18).Value = 4.282548534656
$ws.Cells.Item(7, 19).Value = 0.15052418278737
$ws.Cells.Item(7, 20).Value = 0.15052418278737

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Tnfsf18"
$ws.Cells.Item(8, 3).Value = "Tnfrsf18"
$ws.Cells.Item(8, 4).Value = "MuSCs"
$ws.Cells.Item(8, 5).Value = 2
$ws.Cells.Item(8, 6).Value = 0.6666666666666666
$ws.Cells.Item(8, 7).Value = 0.2127573333333333
$ws.Cells.Item(8, 8).Value = 0.638272
$ws.Cells.Item(8, 9).Value = 0.2734542388551691
$ws.Cells.Item(8, 10).Value = 0.2734542388551692
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 0.998142
$ws.Cells.Item(8, 14).Value = 2.994426
$ws.Cells.Item(8, 15).Value = 0.2456624069636779
$ws.Cells.Item(8, 16).Value = 0.2456624069636779
$ws.Cells.Item(8, 17).Value = 0.212362030208
$ws.Cells.Item(8, 18).Value = 1.911258271872
$ws.Cells.Item(8, 19).Value = 0.06717742651158133
$ws.Cells.Item(8, 20).Value = 0.06717742651158136

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Tnfsf18"
$ws.Cells.Item(9, 3).Value = "Tnfrsf18"
$ws.Cells.Item(9, 4).Value = "Resolving-Mac"
$ws.Cells.Item(9, 5).Value = 2
$ws.Cells.Item(9, 6).Value = 0.6666666666666666
$ws.Cells.Item(9, 7).Value = 0.2127573333333333
$ws.Cells.Item(9, 8).Value = 0.638272
$ws.Cells.Item(9, 9).Value = 0.2734542388551691
$ws.Cells.Item(9, 10).Value = 0.2734542388551692
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.4406919999999999
$ws.Cells.Item(9, 14).Value = 1.322076
$ws.Cells.Item(9, 15).Value = 0.1084629816695792
$ws.Cells.Item(9, 16).Value = 0.1084629816695792
$ws.Cells.Item(9, 17).Value = 0.09376045474133331
$ws.Cells.Item(9, 18).Value = 0.8438440926719998
$ws.Cells.Item(9, 19).Value = 0.02965966209641694
$ws.Cells.Item(9, 20).Value = 0.02965966209641695

# Row 10
$ws.Cells.Item(10, 1).Value = "Resolving-Mac"
$ws.Cells.Item(10, 2).Value = "Tnfsf18"
$ws.Cells.Item(10, 3).Value = "Tnfrsf18"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 2
$ws.Cells.Item(10, 6).Value = 0.6666666666666666
$ws.Cells.Item(10, 7).Value = 0.4408636666666667
$ws.Cells.Item(10, 8).Value = 1.322591
$ws.Cells.Item(10, 9).Value = 0.5666363481739713
$ws.Cells.Item(10, 10).Value = 0.5666363481739713
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 0.6666666666666666
$ws.Cells.Item(10, 13).Value = 0.387697
$ws.Cells.Item(10, 14).Value = 1.163091
$ws.Cells.Item(10, 15).Value = 0.09541986830791312
$ws.Cells.Item(10, 16).Value = 0.09541986830791316
$ws.Cells.Item(10, 17).Value = 0.1709215209756667
$ws.Cells.Item(10, 18).Value = 1.538293688781
$ws.Cells.Item(10, 19).Value = 0.05406836572123715
$ws.Cells.Item(10, 20).Value = 0.05406836572123717

# Row 11
$ws.Cells.Item(11, 1).Value = "Resolving-Mac"
$ws.Cells.Item(11, 2).Value = "Tnfsf18"
$ws.Cells.Item(11, 3).Value = "Tnfrsf18"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 2
$ws.Cells.Item(11, 6).Value = 0.6666666666666666
$ws.Cells.Item(11, 7).Value = 0.4408636666666667
$ws.Cells.Item(11, 8).Value = 1.322591
$ws.Cells.Item(11, 9).Value = 0.5666363481739713
$ws.Cells.Item(11, 10).Value = 0.5666363481739713
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 2.236532666666667
$ws.Cells.Item(11, 14).Value = 6.709598
$ws.Cells.Item(11, 15).Value = 0.5504547430588297
$ws.Cells.Item(11, 16).Value = 0.5504547430588298
$ws.Cells.Item(11, 17).Value = 0.9860059920464446
$ws.Cells.Item(11, 18).Value = 8.874053928418
$ws.Cells.Item(11, 19).Value = 0.311907665441897
$ws.Cells.Item(11, 20).Value = 0.311907665441897

# Row 12
$ws.Cells.Item(12, 1).Value = "Resolving-Mac"
$ws.Cells.Item(12, 2).Value = "Tnfsf18"
$ws.Cells.Item(12, 3).Value = "Tnfrsf18"
$ws.Cells.Item(12, 4).Value = "MuSCs"
$ws.Cells.Item(12, 5).Value = 2
$ws.Cells.Item(12, 6).Value = 0.6666666666666666
$ws.Cells.Item(12, 7).Value = 0.4408636666666667
$ws.Cells.Item(12, 8).Value = 1.322591
$ws.Cells.Item(12, 9).Value = 0.5666363481739713
$ws.Cells.Item(12, 10).Value = 0.5666363481739713
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 0.998142
$ws.Cells.Item(12, 14).Value = 2.994426
$ws.Cells.Item(12, 15).Value = 0.2456624069636779
$ws.Cells.Item(12, 16).Value = 0.2456624069636779
$ws.Cells.Item(12, 17).Value = 0.4400445419740001
$ws.Cells.Item(12, 18).Value = 3.960400877766
$ws.Cells.Item(12, 19).Value = 0.1392012491655264
$ws.Cells.Item(12, 20).Value = 0.1392012491655265

# Row 13
$ws.Cells.Item(13, 1).Value = "Resolving-Mac"
$ws.Cells.Item(13, 2).Value = "Tnfsf18"
$ws.Cells.Item(13, 3).Value = "Tnfrsf18"
$ws.Cells.Item(13, 4).Value = "Resolving-Mac"
$ws.Cells.Item(13, 5).Value = 2
$ws.Cells.Item(13, 6).Value = 0.6666666666666666
$ws.Cells.Item(13, 7).Value = 0.4408636666666667
$ws.Cells.Item(13, 8).Value = 1.322591
$ws.Cells.Item(13, 9).Value = 0.5666363481739713
$ws.Cells.Item(13, 10).Value = 0.5666363481739713
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 0.4406919999999999
$ws.Cells.Item(13, 14).Value = 1.322076
$ws.Cells.Item(13, 15).Value = 0.1084629816695792
$ws.Cells.Item(13, 16).Value = 0.1084629816695792
$ws.Cells.Item(13, 17).Value = 0.1942850909906667
$ws.Cells.Item(13, 18).Value = 1.748565818916
$ws.Cells.Item(13, 19).Value = 0.06145906784531074
$ws.Cells.Item(13, 20).Value = 0.06145906784531076
